$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content in rows 26-46 (A:V) before rewriting with new layout
$ws.Range("A26:V46").ClearContents()

$ws.Range("A26").Value = 'Coroner'
$ws.Range("B26").Value = '検視官'

$ws.Range("A27").Value = 'NoDeadBodies'
$ws.Range("B27").Value = 'この付近に~r~死体~s~はありません。'

$ws.Range("A28").Value = 'UnitRequested'
$ws.Range("B28").Value = '~b~{0}~s~の応援を要請しました。'

$ws.Range("A29").Value = 'CoronerCheckPls'
$ws.Range("B29").Value = '詳しい情報は~b~検視官レポート~s~を確認してください。'

$ws.Range("A30").Value = 'CoronerBye'
$ws.Range("B30").Value = 'それではいい一日を!'

$ws.Range("A31").Value = 'TeleportUnit'
$ws.Range("B31").Value = '{0}で応援を近くにテレポートさせます。'

$ws.Range("A33").Value = 'CoronerMenu'
$ws.Range("B33").Value = '検視官メニュー'

$ws.Range("A34").Value = 'CoronerReport'
$ws.Range("B34").Value = '検視官レポート'

$ws.Range("A35").Value = 'CoronerReportCount'
$ws.Range("B35").Value = '検視官レポート数: {0}'

$ws.Range("A36").Value = 'NoData'
$ws.Range("B36").Value = 'データなし'

$ws.Range("A38").Value = 'Name'
$ws.Range("B38").Value = '名前'

$ws.Range("A39").Value = 'Sex'
$ws.Range("B39").Value = '性別'

$ws.Range("A40").Value = 'CauseOfDeath'
$ws.Range("B40").Value = '死因'

$ws.Range("A41").Value = 'DiedDay'
$ws.Range("B41").Value = '死亡日'

$ws.Range("A43").Value = 'BackupVehicle'
$ws.Range("B43").Value = '応援車両'

$ws.Range("A44").Value = 'BackupOfficer'
$ws.Range("B44").Value = '応援警官'

$ws.Range("A46").Value = 'AllDismissItem'
$ws.Range("B46").Value = '~r~全応援を解散~s~'
$ws.Range("E46").Value = '~r~All Units Dismiss~s~'

